{"js": "// The original paragraph reads:\n//   \"kekler , i\u00e7inde muz ve \u00e7ilek i\u00e7ermelidir.\"\n// The target splits every word/punctuation/space into its own run (the\n// shape Word leaves behind after a spell-check pass), wraps each\n// \"real\" word in <w:proofErr w:type=\"spellStart\"/> ... <w:proofErr\n// w:type=\"spellEnd\"/>, and swaps \"muz ve \u00e7ilek\" (banana and\n// strawberry) for \"\u00e7ikolata ve krema\" (chocolate and cream).\n//\n// Building that exact run/proofErr layout isn't reachable through the\n// higher-level text APIs (insertText/search-replace only ever touch\n// plain runs), so we hand Word the literal OOXML for the new content\n// via insertOoxml \u2014 Word.RequestContext's equivalent of Range.InsertXML.\n// insertOoxml requires the Flat OPC wrapper Office.js always expects.\n\nconst FLAT_OPC_TEMPLATE = (paragraphInnerXml) => `<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>${paragraphInnerXml}</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nfunction spellCheckedWord(word) {\n  return `<w:proofErr w:type=\"spellStart\"/><w:r><w:t>${word}</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>`;\n}\n\nfunction space() {\n  return `<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>`;\n}\n\nfunction plainRun(text) {\n  return `<w:r><w:t>${text}</w:t></w:r>`;\n}\n\nconst newParagraphInnerXml =\n  spellCheckedWord(\"kekler\") +\n  space() +\n  plainRun(\",\") +\n  spellCheckedWord(\"i\u00e7inde\") +\n  space() +\n  spellCheckedWord(\"\u00e7ikolata\") +\n  space() +\n  spellCheckedWord(\"ve\") +\n  space() +\n  spellCheckedWord(\"krema\") +\n  space() +\n  spellCheckedWord(\"i\u00e7ermelidir\") +\n  plainRun(\".\");\n\nconst body = context.document.body;\n\n// Locate the exact original sentence. Searching for the literal text\n// (rather than indexing body.paragraphs) gives us back a Range whose\n// extent is exactly the text run span \u2014 it does NOT include the\n// paragraph mark, so replacing it with insertOoxml leaves the host\n// <w:p> (and its w14:paraId/rsid attributes) untouched.\nconst results = body.search(\"kekler , i\u00e7inde muz ve \u00e7ilek i\u00e7ermelidir.\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to update.\");\n}\n\nresults.items[0].insertOoxml(\n  FLAT_OPC_TEMPLATE(newParagraphInnerXml),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# The original paragraph reads:\n#   \"kekler , i\u00e7inde muz ve \u00e7ilek i\u00e7ermelidir.\"\n# The target splits every word/punctuation/space into its own run (the\n# shape Word leaves behind after a spell-check pass), wraps each \"real\"\n# word in <w:proofErr w:type=\"spellStart\"/> ... <w:proofErr\n# w:type=\"spellEnd\"/>, and swaps \"muz ve \u00e7ilek\" (banana and\n# strawberry) for \"\u00e7ikolata ve krema\" (chocolate and cream).\n#\n# That exact run/proofErr layout can't be produced with Find/Replace or\n# Range.Text (those only ever write a single plain run), so we hand\n# Word the literal OOXML for the new paragraph body via\n# Range.InsertXML, wrapped in the Flat OPC package format Word expects.\n\n$d = $word.ActiveDocument\n\n$newParagraphInnerXml = @'\n<w:proofErr w:type=\"spellStart\"/><w:r><w:t>kekler</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i\u00e7inde</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u00e7ikolata</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ve</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>krema</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i\u00e7ermelidir</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r>\n'@\n\n$flatOpc = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>$newParagraphInnerXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n\"@\n\n# Find the exact original sentence and narrow the range so it excludes\n# the trailing paragraph mark: InsertXML REPLACEs whatever the range\n# spans, and a range that swallowed the pilcrow would blow away the\n# host <w:p> (and its w14:paraId/rsid attributes) along with the text.\n$find = $d.Content.Duplicate\n$find.Find.ClearFormatting()\n$found = $find.Find.Execute(\"kekler , i\u00e7inde muz ve \u00e7ilek i\u00e7ermelidir.\", $true)\n\nif (-not $found) {\n    throw \"Could not find the target sentence to update.\"\n}\n\n$find.InsertXML($flatOpc)\n"}
